# Actualización Automática de Datos (EA1) y (EA2)
# Updates the "timestamp" column (H) for all data rows (2-51) to the new
# refresh timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-03-14 03:54:41"

for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 8).Value = $newTimestamp
}
